$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number by Excel
# are forced to Text format first, then reset back to the default style so no
# stray number formatting is left behind on the cell.

$ws.Range("D2").Value = "28.057.03"
$ws.Range("E2").Value = "  -3.57%  "
$ws.Range("D3").Value = "1.747.50"
$ws.Range("E3").Value = "  -4.17%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5805"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2709"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06585"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07511"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").Value = "1.748.11"
$ws.Range("E12").Value = "  -4.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.709"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6021"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.60%  "
$ws.Range("D15").Value = "1.986.14"
$ws.Range("E15").Value = "  -4.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "73.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008640"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -11.84%  "
$ws.Range("D18").Value = "28.079.33"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.315"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "204.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.665"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.007"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1233"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.08%  "
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.384"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06113"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.386"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.728"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.707"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.672"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.036"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6334"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.447"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.709"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01679"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.283"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.65%  "
$ws.Range("D41").Value = "1.125.39"
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8655"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.003"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("D45").Value = "1.898.93"
$ws.Range("E45").Value = "  -3.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.575"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000107"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.282"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05395"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4435"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.12%  "
